$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level window position (best effort) ---
$wb.Windows.Item(1).Left = 760
$wb.Windows.Item(1).Top = 760

# --- Remove the duplicated hidden chart helper names (v1.2/v1.3 were
#     accidental dupes of v1.0/v1.1) ---
foreach ($dupName in @("_xlchart.v1.2", "_xlchart.v1.3")) {
    try {
        $wb.Names.Item($dupName).Delete()
    } catch {
    }
}

# --- New "increase" header labels (row 18), bold like the other headers ---
$ws.Range("D18").Value = "Mean increase"
$ws.Range("D18").Font.Bold = $true
$ws.Range("F18").Value = "Median increase"
$ws.Range("F18").Font.Bold = $true

# --- New formulas (row 19): percentage increase vs. a reference baseline ---
$ws.Range("D19").Formula = "=((E3/114.202998)*100)-100"
$ws.Range("D19").ClearFormats()
$ws.Range("F19").Formula = "=((E10/113.658804)*100)-100"
$ws.Range("F19").ClearFormats()

# --- Update the selected cell to match the author's final cursor position ---
$ws.Range("E24").Select() | Out-Null
